# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price column to Text format so that numeric-looking
# strings (e.g. "313.65") are stored as text instead of being coerced to
# floating point numbers, matching the inline-string cells already in the sheet.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.236.07"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.882.59"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "313.65"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.5135"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").Value = "0.3905"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("D9").Value = "0.08381"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").Value = "41.59"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "6.245"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "1.884.70"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "7.306"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "91.47"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "0.06664"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "6.062"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").Value = "28.271.55"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "2.270"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "2.092.73"
$ws.Range("D27").Value = "2.517"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").Value = "158.96"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "20.67"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "125.54"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").Value = "1.043"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "5.896"
$ws.Range("D34").Value = "3.594"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "9.787"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "0.02463"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "0.06564"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "0.2195"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "0.6529"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "5.038"
$ws.Range("E41").Value = "  +3.27%  "
$ws.Range("D42").Value = "1.232"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "0.6126"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "13.16"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "1.293"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "3.682"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "2.024"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "1.233"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "78.81"
$ws.Range("E51").Value = "  -1.61%  "

# Restore the original (default) cell style now that the text values are set.
$priceRange.Style = "Normal"
